# Generate Report for Handback
# Applies the "handback" localization-status update:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#    columns for the zh-cn and de-de target-language tabs, including hyperlinks on the new
#    "Latest Target File" cells
#  - Widens a handful of columns that now hold longer text

$wb = $excel.ActiveWorkbook

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/644a445072deff8aff57fa09e1cd19b7e7267779/e2e/686519fd-cb20-4ac0-a26f-319f30ea898d.md"
$targetMdName = "686519fd-cb20-4ac0-a26f-319f30ea898d.md"

# ---------------------------------------------------------------------------
# 1. Overview sheet: status text + column widths
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($addr in "E2", "F2", "E3", "F3") {
    $wsOverview.Range($addr).Value = "Handed back: in sync with en-US"
}

$wsOverview.Range("E1").ColumnWidth = 29.144371396019366
$wsOverview.Range("F1").ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File (I) + hyperlink
$wsZhCn.Range("I2").Value = $targetMdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetMdUrl, "", "", $targetMdName) | Out-Null
$wsZhCn.Range("I3").Value = $targetMdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $targetMdUrl, "", "", $targetMdName) | Out-Null

# Latest Handback File (J)
$wsZhCn.Range("J2").Value = "686519fd-cb20-4ac0-a26f-319f30ea898d.d074816ff4c5ba13c4643940ba51c2fe03a3eeae.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "686519fd-cb20-4ac0-a26f-319f30ea898d.d074816ff4c5ba13c4643940ba51c2fe03a3eeae.zh-cn.xlf"

# Latest Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-09-01 15:28:52"
$wsZhCn.Range("K3").Value = "2016-09-01 15:28:52"

# Column widths
$wsZhCn.Range("C1").ColumnWidth = 29.144371396019366
$wsZhCn.Range("I1").ColumnWidth = 39.166666666666664
$wsZhCn.Range("J1").ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Target File (I) + hyperlink
$wsDeDe.Range("I2").Value = $targetMdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetMdUrl, "", "", $targetMdName) | Out-Null
$wsDeDe.Range("I3").Value = $targetMdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $targetMdUrl, "", "", $targetMdName) | Out-Null

# Latest Handback File (J)
$wsDeDe.Range("J2").Value = "686519fd-cb20-4ac0-a26f-319f30ea898d.d074816ff4c5ba13c4643940ba51c2fe03a3eeae.de-de.xlf"
$wsDeDe.Range("J3").Value = "686519fd-cb20-4ac0-a26f-319f30ea898d.d074816ff4c5ba13c4643940ba51c2fe03a3eeae.de-de.xlf"

# Latest Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-09-01 15:29:00"
$wsDeDe.Range("K3").Value = "2016-09-01 15:29:00"

# Column widths
$wsDeDe.Range("C1").ColumnWidth = 29.144371396019366
$wsDeDe.Range("I1").ColumnWidth = 39.166666666666664
$wsDeDe.Range("J1").ColumnWidth = 39.166666666666664

$wb.Save()
